$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing odds values (rows 2, 4, 6, 8, 11) ---
$ws.Range("O2").Value = 1.62
$ws.Range("P2").Value = 2.2
$ws.Range("G4").Value = 1.6
$ws.Range("N4").Value = 7.5
$ws.Range("AC4").Value = 7.5
$ws.Range("AD4").Value = 7
$ws.Range("AO4").Value = 8.5
$ws.Range("AQ4").Value = 29
$ws.Range("AU4").Value = 9.5
$ws.Range("G6").Value = 1.91
$ws.Range("I6").Value = 3.7
$ws.Range("J6").Value = 2.5
$ws.Range("L6").Value = 4
$ws.Range("Q6").Value = 1.67
$ws.Range("R6").Value = 2.15
$ws.Range("W6").Value = 9
$ws.Range("X6").Value = 10
$ws.Range("Y6").Value = 8.5
$ws.Range("Z6").Value = 17
$ws.Range("AE6").Value = 13
$ws.Range("AJ6").Value = 13
$ws.Range("AN6").Value = 4
$ws.Range("AO6").Value = 10
$ws.Range("G8").Value = 1.95
$ws.Range("I8").Value = 4
$ws.Range("J8").Value = 2.63
$ws.Range("M8").Value = 1.06
$ws.Range("N8").Value = 10
$ws.Range("X8").Value = 9
$ws.Range("Z8").Value = 17
$ws.Range("AI8").Value = 19
$ws.Range("AJ8").Value = 13
$ws.Range("AN8").Value = 4
$ws.Range("AO8").Value = 11
$ws.Range("AP8").Value = 23
$ws.Range("AQ8").Value = 41
$ws.Range("AW8").Value = 5.5
$ws.Range("AX8").Value = 21
$ws.Range("AY8").Value = 29
$ws.Range("G11").Value = 3
$ws.Range("I11").Value = 2.55
$ws.Range("L11").Value = 3.4
$ws.Range("M11").Value = 1.1
$ws.Range("N11").Value = 7
$ws.Range("Z11").Value = 34
$ws.Range("AW11").Value = 4.33
$ws.Range("AX11").Value = 15
$ws.Range("BA11").Value = 81

# --- Add new row 13 (Racing Montevideo - Danubio) ---
$ws.Range("A13").Value = "0Mk66xgo"
$ws.Range("B13").Value = "17/11/2024"
$ws.Range("C13").Value = "16:30"
$ws.Range("D13").Value = "URUGUAY - PRIMERA DIVISION"
$ws.Range("E13").Value = "Racing Montevideo"
$ws.Range("F13").Value = "Danubio"
$ws.Range("G13").Value = 2.55
$ws.Range("H13").Value = 2.9
$ws.Range("I13").Value = 3
$ws.Range("J13").Value = 3.5
$ws.Range("K13").Value = 1.83
$ws.Range("L13").Value = 4
$ws.Range("M13").Value = 1.13
$ws.Range("N13").Value = 6
$ws.Range("O13").Value = 1.53
$ws.Range("P13").Value = 2.38
$ws.Range("Q13").Value = 2.7
$ws.Range("R13").Value = 1.44
$ws.Range("S13").Value = 1.62
$ws.Range("T13").Value = 2.2
$ws.Range("U13").Value = 2.2
$ws.Range("V13").Value = 1.62
$ws.Range("W13").Value = 6
$ws.Range("X13").Value = 11
$ws.Range("Y13").Value = 11
$ws.Range("Z13").Value = 26
$ws.Range("AA13").Value = 26
$ws.Range("AB13").Value = 41
$ws.Range("AC13").Value = 6
$ws.Range("AD13").Value = 6
$ws.Range("AE13").Value = 19
$ws.Range("AF13").Value = 81
$ws.Range("AG13").Value = 1250
$ws.Range("AH13").Value = 7
$ws.Range("AI13").Value = 13
$ws.Range("AJ13").Value = 12
$ws.Range("AK13").Value = 34
$ws.Range("AL13").Value = 29
$ws.Range("AM13").Value = 41
$ws.Range("AN13").Value = 4.33
$ws.Range("AO13").Value = 17
$ws.Range("AP13").Value = 34
$ws.Range("AQ13").Value = 51
$ws.Range("AR13").Value = 101
$ws.Range("AS13").Value = 351
$ws.Range("AT13").Value = 2.2
$ws.Range("AU13").Value = 9.5
$ws.Range("AV13").Value = 81
$ws.Range("AW13").Value = 4.75
$ws.Range("AX13").Value = 19
$ws.Range("AY13").Value = 34
$ws.Range("AZ13").Value = 67
$ws.Range("BA13").Value = 101
$ws.Range("BB13").Value = 351
$ws.Range("BC13").Value = 51
$ws.Range("BD13").Value = 51
